$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from "joanne" to "praveen"
$ws.Range("B2").Value = "praveen"

# Update B6 from "praveen" to "lol"
$ws.Range("B6").Value = "lol"

# Add new row 7
$ws.Range("A7").Value = "lol"
$ws.Range("B7").Value = "lol"
$ws.Range("C7").Value = "l"
$ws.Range("D7").Value = "2l"
$ws.Range("E7").Value = "3l"
$ws.Range("F7").Value = "4l"
$ws.Range("G7").Value = "5l"

# Update selection to G8
$ws.Range("G8").Select()
